$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "PN"
$ws.Range("B3").Value = "CF"
$ws.Range("B5").Value = "UR"
$ws.Range("B6").Value = "0O"
$ws.Range("B7").Value = "NGJ"
$ws.Range("B8").Value = "A8"
$ws.Range("B13").Value = "0Z"
